{"js": "// Replace the 25 \"AxB=C\" answer strings in the multiplication table with\n// their new values, as described by the diff. Each old value is unique in\n// the document, so a straightforward exact-text search & replace for every\n// pair is safe (no ambiguity / no cross-matches between old and new values).\nconst replacements = [\n  [\"66\u00d792=6072\", \"23\u00d740=920\"],\n  [\"36\u00d751=1836\", \"30\u00d733=990\"],\n  [\"18\u00d740=720\", \"98\u00d783=8134\"],\n  [\"91\u00d731=2821\", \"22\u00d728=616\"],\n  [\"53\u00d716=848\", \"91\u00d791=8281\"],\n  [\"11\u00d722=242\", \"14\u00d772=1008\"],\n  [\"49\u00d727=1323\", \"38\u00d762=2356\"],\n  [\"72\u00d786=6192\", \"76\u00d765=4940\"],\n  [\"59\u00d757=3363\", \"87\u00d719=1653\"],\n  [\"71\u00d740=2840\", \"49\u00d747=2303\"],\n  [\"20\u00d726=520\", \"72\u00d788=6336\"],\n  [\"50\u00d756=2800\", \"40\u00d720=800\"],\n  [\"54\u00d771=3834\", \"26\u00d730=780\"],\n  [\"91\u00d793=8463\", \"48\u00d738=1824\"],\n  [\"96\u00d789=8544\", \"15\u00d714=210\"],\n  [\"44\u00d787=3828\", \"17\u00d782=1394\"],\n  [\"38\u00d732=1216\", \"89\u00d764=5696\"],\n  [\"74\u00d789=6586\", \"18\u00d743=774\"],\n  [\"55\u00d773=4015\", \"41\u00d713=533\"],\n  [\"88\u00d725=2200\", \"40\u00d790=3600\"],\n  [\"20\u00d790=1800\", \"97\u00d769=6693\"],\n  [\"15\u00d752=780\", \"73\u00d722=1606\"],\n  [\"37\u00d766=2442\", \"43\u00d781=3483\"],\n  [\"74\u00d731=2294\", \"33\u00d751=1683\"],\n  [\"97\u00d799=9603\", \"36\u00d748=1728\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"AxB=C\" answer strings in the multiplication table with\n# their new values, as described by the diff. Each old value is unique in\n# the document, so an exact Find/Replace-All pass per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"66\u00d792=6072\", \"23\u00d740=920\"),\n    @(\"36\u00d751=1836\", \"30\u00d733=990\"),\n    @(\"18\u00d740=720\", \"98\u00d783=8134\"),\n    @(\"91\u00d731=2821\", \"22\u00d728=616\"),\n    @(\"53\u00d716=848\", \"91\u00d791=8281\"),\n    @(\"11\u00d722=242\", \"14\u00d772=1008\"),\n    @(\"49\u00d727=1323\", \"38\u00d762=2356\"),\n    @(\"72\u00d786=6192\", \"76\u00d765=4940\"),\n    @(\"59\u00d757=3363\", \"87\u00d719=1653\"),\n    @(\"71\u00d740=2840\", \"49\u00d747=2303\"),\n    @(\"20\u00d726=520\", \"72\u00d788=6336\"),\n    @(\"50\u00d756=2800\", \"40\u00d720=800\"),\n    @(\"54\u00d771=3834\", \"26\u00d730=780\"),\n    @(\"91\u00d793=8463\", \"48\u00d738=1824\"),\n    @(\"96\u00d789=8544\", \"15\u00d714=210\"),\n    @(\"44\u00d787=3828\", \"17\u00d782=1394\"),\n    @(\"38\u00d732=1216\", \"89\u00d764=5696\"),\n    @(\"74\u00d789=6586\", \"18\u00d743=774\"),\n    @(\"55\u00d773=4015\", \"41\u00d713=533\"),\n    @(\"88\u00d725=2200\", \"40\u00d790=3600\"),\n    @(\"20\u00d790=1800\", \"97\u00d769=6693\"),\n    @(\"15\u00d752=780\", \"73\u00d722=1606\"),\n    @(\"37\u00d766=2442\", \"43\u00d781=3483\"),\n    @(\"74\u00d731=2294\", \"33\u00d751=1683\"),\n    @(\"97\u00d799=9603\", \"36\u00d748=1728\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
